$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

$rows = @(
    @(18, 15, "NONE", "S",    "O",      "Note"),
    @(19, 16, "NONE", "047G", "047G-P", "Huomenna on RYÖpäivä"),
    @(20, 17, "NONE", "S",    "O",      "Best taustakuva"),
    @(21, 18, "NONE", "S",    "O",      "Note"),
    @(22, 19, "NONE", "S",    "O",      "Note"),
    @(23, 20, "NONE", "S",    "O",      "Note")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
